# Atualizei dados da bibi e add
# Applies the refreshed "previsao_retorno" figures to the Resumo_por_Cliente sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "situacao" (column J) text refreshed for several clients whose
#     inactivity window (meses sem comprar) advanced slightly ---
$ws.Range("J4").Value   = "INATIVO - 36.5 meses sem comprar"
$ws.Range("J9").Value   = "INATIVO - 19.5 meses sem comprar"
$ws.Range("J50").Value  = "INATIVO - 7.2 meses sem comprar"
$ws.Range("J52").Value  = "INATIVO - 10.4 meses sem comprar"
$ws.Range("J67").Value  = "INATIVO - 22.3 meses sem comprar"
$ws.Range("J69").Value  = "INATIVO - 13.0 meses sem comprar"
$ws.Range("J75").Value  = "INATIVO - 22.3 meses sem comprar"
$ws.Range("J83").Value  = "INATIVO - 22.7 meses sem comprar"
$ws.Range("J87").Value  = "INATIVO - 22.2 meses sem comprar"
$ws.Range("J88").Value  = "INATIVO - 9.9 meses sem comprar"
$ws.Range("J93").Value  = "INATIVO - 15.5 meses sem comprar"
$ws.Range("J98").Value  = "INATIVO - 19.6 meses sem comprar"
$ws.Range("J102").Value = "INATIVO - 23.3 meses sem comprar"
$ws.Range("J111").Value = "INATIVO - 21.8 meses sem comprar"
$ws.Range("J112").Value = "INATIVO - 7.5 meses sem comprar"

# --- Row 38 (id_cliente 4165, JOAO VITOR MARQUES REIS) bought again,
#     so its purchase-probability range, purchase count and the
#     ultima/proxima compra window all moved forward ---
$ws.Range("B38").Value = 0.33
$ws.Range("C38").Value = 0.17
$ws.Range("E38").Value = 22
$ws.Range("H38").Value = 45848.72115740741
$ws.Range("I38").Value = 45910.72115740741

# --- Row 115 (id_cliente 28458, BEMOL S/A) also has an updated
#     purchase history count and refreshed ultima/proxima compra ---
$ws.Range("E115").Value = 16515
$ws.Range("H115").Value = 45848.67685185185
$ws.Range("I115").Value = 45849.67685185185
